# Split the three concatenated "Programa"/"Bibliografia" paragraphs into
# multiple <w:t> runs separated by manual line breaks (<w:br/>), matching
# the numbered-list / reference-list formatting applied upstream.
$d = $word.ActiveDocument

$found1 = $d.Content.Find.Execute(
    "1) Introdução a Física: sistemas de unidades, revisão de vetores, análise dimensional.2) Cinemática: movimento unidimensional, queda livre, movimento bidimensional, projéteis. 3) Dinâmica: leis de Newton, forças, força de atrito, força de resistência do ar, velocidade terminal, movimento circular uniforme, gravitação, aplicações.4) Energia: trabalho, forças conservativas, conservação de energia mecânica, atrito, aplicações.5)  Momento linear: centro de massa, sistema de partículas, conservação do momento linear, colisões, impulso.6) Rotação: variáveis do movimento rotacional, energia cinética rotacional, momento de inércia, torque, rolamento, conservação do momento angular.",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "1) Introdução a Física: sistemas de unidades, revisão de vetores, análise dimensional.^l2) Cinemática: movimento unidimensional, queda livre, movimento bidimensional, projéteis. ^l3) Dinâmica: leis de Newton, forças, força de atrito, força de resistência do ar, velocidade terminal, movimento circular uniforme, gravitação, aplicações.^l4) Energia: trabalho, forças conservativas, conservação de energia mecânica, atrito, aplicações.^l5)  Momento linear: centro de massa, sistema de partículas, conservação do momento linear, colisões, impulso.^l6) Rotação: variáveis do movimento rotacional, energia cinética rotacional, momento de inércia, torque, rolamento, conservação do momento angular.", 2)
if (-not $found1) { throw "PT Programa paragraph not found" }
Write-Host "PT Programa paragraph: line breaks inserted =" $found1

$found2 = $d.Content.Find.Execute(
    "1) Introduction to Physics: unit systems, review of vectors, dimensional analysis. 2) Kinematics: one dimensional motion, free fall, bidimensional motion, projectile.  3) Dynamics: Newton’s laws, friction force, drag force, terminal speed, uniform circular motion, gravitation, applications.4) Energy: work, conservative forces, mechanical energy conservation, friction, applications.5)  Linear momentum: center of mass, system of particles, conservation of linear momentum, collisions, impulse.6) Rotation: rotational variables, kinetic energy of rotation, rotational inertia, torque, rolling, conservation of angular momentum",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "1) Introduction to Physics: unit systems, review of vectors, dimensional analysis. ^l2) Kinematics: one dimensional motion, free fall, bidimensional motion, projectile. ^l 3) Dynamics: Newton’s laws, friction force, drag force, terminal speed, uniform circular motion, gravitation, applications.^l4) Energy: work, conservative forces, mechanical energy conservation, friction, applications.^l5)  Linear momentum: center of mass, system of particles, conservation of linear momentum, collisions, impulse.^l6) Rotation: rotational variables, kinetic energy of rotation, rotational inertia, torque, rolling, conservation of angular momentum", 2)
if (-not $found2) { throw "EN Programa paragraph not found" }
Write-Host "EN Programa paragraph: line breaks inserted =" $found2

$found3 = $d.Content.Find.Execute(
    "HALLIDAY, D; RESNICK, R. Fundamentos de Física. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).",
    $true, $true, $false, $false, $false, $true, 1, $false,
    "HALLIDAY, D; RESNICK, R. Fundamentos de Física. Vol.1, LTC (2008).^lSEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).^lJEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).^lNUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).^lTIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).", 2)
if (-not $found3) { throw "Bibliografia paragraph not found" }
Write-Host "Bibliografia paragraph: line breaks inserted =" $found3
